$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1275097979829652
$ws.Range("C2").Value = 1.514122279856028
$ws.Range("D2").Value = 9.324398268036724
$ws.Range("E2").Value = 3.05358776982695
$ws.Range("F2").Value = 3.122720650086514
$ws.Range("G2").Value = 22

$ws.Range("B3").Value = 0.029890368269258
$ws.Range("C3").Value = 1.46966240847316
$ws.Range("D3").Value = 8.914287279745698
$ws.Range("E3").Value = 2.985680371330075
$ws.Range("F3").Value = 3.059258658222928
$ws.Range("G3").Value = 21

$ws.Range("B4").Value = -0.4518501121711197
$ws.Range("C4").Value = 0.9211732945268647
$ws.Range("D4").Value = 3.928571433859458
$ws.Range("E4").Value = 1.982062419264201
$ws.Range("F4").Value = 1.980006089833911
$ws.Range("G4").Value = 20

$ws.Range("B5").Value = -0.02887321458375646
$ws.Range("C5").Value = 0.6406565445890792
$ws.Range("D5").Value = 0.9487813136001403
$ws.Range("E5").Value = 0.9740540609227705
$ws.Range("F5").Value = 1.000305658023116
$ws.Range("G5").Value = 19

$ws.Range("B6").Value = -0.01274653667201224
$ws.Range("C6").Value = 0.6862274353315703
$ws.Range("D6").Value = 0.9980463954621034
$ws.Range("E6").Value = 0.9990227201931412
$ws.Range("F6").Value = 1.027902220766659
$ws.Range("G6").Value = 18

$ws.Range("B7").Value = -0.1249833021668538
$ws.Range("C7").Value = 0.5228445079367579
$ws.Range("D7").Value = 0.5102717666494121
$ws.Range("E7").Value = 0.7143330922261771
$ws.Range("F7").Value = 0.7249597400067718
$ws.Range("G7").Value = 17

$ws.Range("B8").Value = -0.0390809901440984
$ws.Range("C8").Value = 0.4768323812870599
$ws.Range("D8").Value = 0.4090940898259321
$ws.Range("E8").Value = 0.6396046355569447
$ws.Range("F8").Value = 0.6593465581702652
$ws.Range("G8").Value = 16
